# Update "想去人数" (col F) and "最低票价" (col G) figures on the
# "展览" and "全部类型" sheets to the freshly scraped values.

$wb = $excel.ActiveWorkbook

# Row => @(new F value, new G value)  -- $null means "leave unchanged"
$updates = @{
    2  = @(280,  55)
    3  = @(1436, 65)
    6  = @(240,  $null)
    12 = @(4771, $null)
    14 = @(7039, $null)
    18 = @(582,  $null)
    20 = @(5,    $null)
    21 = @(4185, $null)
    22 = @(1188, $null)
    23 = @(83,   45)
    25 = @(2767, $null)
    30 = @(393,  $null)
    31 = @(413,  $null)
    34 = @(1653, $null)
    35 = @(1064, $null)
    37 = @(796,  $null)
    45 = @(419,  $null)
    47 = @(25,   $null)
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $vals = $updates[$row]
        $fVal = $vals[0]
        $gVal = $vals[1]
        if ($null -ne $fVal) {
            $ws.Cells.Item($row, 6).Value = $fVal
        }
        if ($null -ne $gVal) {
            $ws.Cells.Item($row, 7).Value = $gVal
        }
    }
}
